# Update the "想去人数" (F) and "最低票价" (G) columns on the
# "展览" and "全部类型" sheets, which share identical data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # F column (想去人数) numeric updates
    $ws.Range("F5").Value()  = 471
    $ws.Range("F9").Value()  = 129
    $ws.Range("F12").Value() = 340
    $ws.Range("F13").Value() = 1821
    $ws.Range("F17").Value() = 714
    $ws.Range("F19").Value() = 350
    $ws.Range("F20").Value() = 4360
    $ws.Range("F21").Value() = 19
    $ws.Range("F22").Value() = 315
    $ws.Range("F23").Value() = 1178
    $ws.Range("F24").Value() = 509
    $ws.Range("F26").Value() = 729
    $ws.Range("F28").Value() = 379

    # G column (最低票价) updates
    $ws.Range("G6").Value()  = 55
    $ws.Range("G7").Value()  = "已售罄"
    $ws.Range("G8").Value()  = 65
}
